# Dr. E version from 2019-06-04
# Adds two new "raw -> normalized" type mappings to the ent_type.csv sheet:
#   - "academic sector"  -> "academic institution"   (inserted in sorted position, row 8)
#   - "un entity"         -> "united nations entity"  (appended at the end, row 31)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the "un entity" row at the end of the table first (row 30, right after the
#     existing last row 29) - this is typed before "academic sector" below, which keeps
#     the shared-string table append order matching the authoring sequence. After the
#     later row-8 insert shifts everything down by one, this becomes row 31. ---
$ws.Range("A30").Value = "un entity"
$ws.Range("B30").Value = "united nations entity"

# Match the formatting used elsewhere in that block (left-aligned Arial 10, same as
# the rest of the "type"/"new type" columns) by copying the format from a sibling row.
$ws.Range("B10").Copy()
$ws.Range("A30:B30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Insert the "academic sector" row right before the existing "business" row ---
$ws.Rows(8).Insert()
$ws.Range("A8").Value = "academic sector"
$ws.Range("B8").Value = "academic institution"

# --- Restore the selection to where the editor was working ---
$ws.Range("A9").Select() | Out-Null
